$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34: reference doc path (entered first) ---
$ws.Range("B34").Value = "D:\R\git\all-languages\Unit_Test_Frameworks\Microsoft_Fakes\Choose between stub and shim types.docx"

# --- Row 21: highlighted marker (reuse existing yellow-fill style) + heading ---
$ws.Range("A21").Interior.Color = 65535
$ws.Range("B21").Value = "Choose between stub and shim types"

# --- Row 19: follow-up note under the existing "inheritance" bullet ---
$ws.Range("B19").Value = "More details in the next point"

# --- Row 23-26: STUBS section ---
$ws.Range("B23").Value = "STUBS"
$ws.Range("C24").Value = "Better to have stubs for all the classes in your solution because you have control over them and can implement interfaces"
$ws.Range("C25").Value = "STUBS only work with interfaces"
$ws.Range("C26").Value = "Can provide alternative implementations for the members defined in this interface."

# --- Row 28: SHIMS section heading ---
$ws.Range("B28").Value = "SHIMS"

# --- Row 33: reference sample path ---
$ws.Range("B33").Value = "D:\R\git\all-languages\Unit_Test_Frameworks\Microsoft_Fakes\Samples\Stub_Method_Property\ServicesTests\Images_From_ObjectBrowser"

# --- Row 29: rich-text note about shims ---
$rt = $ws.Range("C29")
$rt.Value = "external assemblies such as System.dll typically are not provided with separate interface definitions, so you must use shims instead."
$rt.Font.Size = 10
$rt.Font.Color = 1513239
$rt.Font.Name = "Segoe UI"
$rt.Characters(29, 10).Font.Italic = $true
$rt.Characters(39, 95).Font.Name = "Segoe UI"
$ws.Rows(29).RowHeight = 15

# --- Row 30 ---
$ws.Range("C30").Value = "If you don't have control over the code in which you can't implement interfaces then its better to go with SHIMS in this cases only."

# --- Row 32 ---
$ws.Range("B32").Value = "FOR MORE INORMATION REFER THESE PATHS IN GIT"

# --- Page setup (printable area attributes) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- View state: scroll + final selection ---
$ws.Activate()
$ws.Range("A36").Select()
